$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 32,4
$data[0,0] = 32
$data[0,1] = "house/house020.jpg"
$data[0,2] = "heißen"
$data[0,3] = "house"
$data[1,0] = 76
$data[1,1] = "dog/dog026.jpg"
$data[1,2] = "tollen"
$data[1,3] = "dog"
$data[2,0] = 19
$data[2,1] = "house/house019.jpg"
$data[2,2] = "bauen"
$data[2,3] = "house"
$data[3,0] = 31
$data[3,1] = "house/house006.jpg"
$data[3,2] = "kriegen"
$data[3,3] = "house"
$data[4,0] = 12
$data[4,1] = "house/house030.jpg"
$data[4,2] = "bergen"
$data[4,3] = "house"
$data[5,0] = 111
$data[5,1] = "dog/dog004.jpg"
$data[5,2] = "lernen"
$data[5,3] = "dog"
$data[6,0] = 43
$data[6,1] = "dog/dog015.jpg"
$data[6,2] = "legen"
$data[6,3] = "dog"
$data[7,0] = 39
$data[7,1] = "dog/dog010.jpg"
$data[7,2] = "danken"
$data[7,3] = "dog"
$data[8,0] = 48
$data[8,1] = "house/house023.jpg"
$data[8,2] = "lassen"
$data[8,3] = "house"
$data[9,0] = 106
$data[9,1] = "house/house005.jpg"
$data[9,2] = "rechnen"
$data[9,3] = "house"
$data[10,0] = 30
$data[10,1] = "house/house013.jpg"
$data[10,2] = "küssen"
$data[10,3] = "house"
$data[11,0] = 100
$data[11,1] = "dog/dog031.jpg"
$data[11,2] = "nullen"
$data[11,3] = "dog"
$data[12,0] = 56
$data[12,1] = "house/house011.jpg"
$data[12,2] = "süßen"
$data[12,3] = "house"
$data[13,0] = 99
$data[13,1] = "dog/dog012.jpg"
$data[13,2] = "wecken"
$data[13,3] = "dog"
$data[14,0] = 34
$data[14,1] = "dog/dog000.jpg"
$data[14,2] = "passen"
$data[14,3] = "dog"
$data[15,0] = 11
$data[15,1] = "dog/dog024.jpg"
$data[15,2] = "hassen"
$data[15,3] = "dog"
$data[16,0] = 109
$data[16,1] = "dog/dog016.jpg"
$data[16,2] = "meinen"
$data[16,3] = "dog"
$data[17,0] = 84
$data[17,1] = "house/house029.jpg"
$data[17,2] = "trotzen"
$data[17,3] = "house"
$data[18,0] = 121
$data[18,1] = "dog/dog028.jpg"
$data[18,2] = "parken"
$data[18,3] = "dog"
$data[19,0] = 0
$data[19,1] = "house/house014.jpg"
$data[19,2] = "heben"
$data[19,3] = "house"
$data[20,0] = 61
$data[20,1] = "house/house018.jpg"
$data[20,2] = "öffnen"
$data[20,3] = "house"
$data[21,0] = 46
$data[21,1] = "house/house017.jpg"
$data[21,2] = "ändern"
$data[21,3] = "house"
$data[22,0] = 13
$data[22,1] = "house/house024.jpg"
$data[22,2] = "deuten"
$data[22,3] = "house"
$data[23,0] = 117
$data[23,1] = "dog/dog022.jpg"
$data[23,2] = "kennen"
$data[23,3] = "dog"
$data[24,0] = 78
$data[24,1] = "dog/dog025.jpg"
$data[24,2] = "achten"
$data[24,3] = "dog"
$data[25,0] = 10
$data[25,1] = "dog/dog017.jpg"
$data[25,2] = "mögen"
$data[25,3] = "dog"
$data[26,0] = 72
$data[26,1] = "dog/dog027.jpg"
$data[26,2] = "sparen"
$data[26,3] = "dog"
$data[27,0] = 70
$data[27,1] = "house/house027.jpg"
$data[27,2] = "hacken"
$data[27,3] = "house"
$data[28,0] = 63
$data[28,1] = "house/house026.jpg"
$data[28,2] = "atmen"
$data[28,3] = "house"
$data[29,0] = 113
$data[29,1] = "house/house007.jpg"
$data[29,2] = "dienen"
$data[29,3] = "house"
$data[30,0] = 105
$data[30,1] = "dog/dog029.jpg"
$data[30,2] = "betteln"
$data[30,3] = "dog"
$data[31,0] = 15
$data[31,1] = "dog/dog013.jpg"
$data[31,2] = "wachsen"
$data[31,3] = "dog"

$ws.Range("B2:E33").Value = $data
